$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New notes / comments (new shared strings + cells) ---
$ws.Range("H1").Value = "This example / end to end test shows the following"
$ws.Range("I1").Value = "Nested Property Setup"
$ws.Range("I2").Value = "Table property setup (using ""table of"")"
$ws.Range("I3").Value = "Nested assertions"
$ws.Range("H5").Value = "As long as there is a gap after the property columns, you can use the rest of the sheet to add notes and do calculations"

# --- Column widths for the new note columns ---
$ws.Columns("H").ColumnWidth = 32.5703125
$ws.Columns("I").ColumnWidth = 32.7109375

# --- Selection, matching the target workbook view ---
$ws.Range("A19:C20").Select
